$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update simulation result values in column C (slight re-measurement) ---
$ws.Range("C3").Value = 0.899997
$ws.Range("C4").Value = 0.100003
$ws.Range("C5").Value = 0.100003
$ws.Range("C6").Value = 0.2
$ws.Range("C7").Value = 0.100003
$ws.Range("C8").Value = 0.100016

# --- Add "Expected results" column D ---
$ws.Range("D3").Value = 0.9
$ws.Range("D4").Value = 0.1
$ws.Range("D5").Value = 0.1
$ws.Range("D6").Value = 0.2
$ws.Range("D7").Value = 0.1
$ws.Range("D8").Value = 0.1
$ws.Range("D3:D8").NumberFormat = "0.00%"

# --- Add "Comment" column E with explanation for each predictor ---
$ws.Range("E3").Value = "Main loop is always taken and it’s 1/10 of all mispredictions => 10% of 90% is main loop mispredicts. Tested branch is for() with 9 stages and one of them is NT. => 1/10(of all predictions) is predicted and 8/10(of all predictions) is mispredicted. => 1/10 + 8/10 = 90%"
$ws.Range("E4").Value = "Main loop is always taken and it’s 1/10 of all mispredictions => 20% of 80% is main loop mispredicts. Tested branch is for() with 4 stages and one of them is NT. => 1/5(of all predictions) is predicted and 3/5(of all predictions) is mispredicted. => 3/5 + 1/5 = 80%"
$ws.Range("E5").Value = "Main loop is always taken backward jump(unconditional) and it’s 1/10 of all predictions => 10% is right. Tested branch is for() with backward jump and it consists from 9 stages and one of them is NT. => 8/9 is predicted and 1/10 is mispredicted. 1/10 = 10%"
$ws.Range("E6").Value = "Main loop is always taken and predictor knows it after the first time => 1/10 predictions is right. Tested branch is jump from for() and has 9 stages: T, T, T, T, T, T, T, T, NT => NT → T → T→ T→ T → T→ T→ T →  NT… We have 7 right predicts and 2 fails(T when NT and after this NT when T). => only 2/10 is mispredicted  and other 6/10 is right predicted. => 2/10 = 20%"
$ws.Range("E7").Value = "Main loop is always taken and predictor knows it after the first time => 1/10 predictions is right. Tested branch is jump from for() and has 9 stages: T, T, T, T, T, T, T, T, NT => WNT → WT → ST→ ST→ ST→ ST→ST→ ST→ST → WT → ST → ST… We have 8 right predicts and 1 fail(ST when NT). => only 1/10 is mispredicted and other 8/10 is right predicted. => 1/10 = 10%"
$ws.Range("E8").Value = "Main loop is always taken and predictor knows it after the first two times => 1/10 of all predictions is right(0 misses). Tested branch is for() with 9 conditions: 8 T and 1 NT => We have history 11111111011111111011…. And we use 2 last bit history. But pattern 11 has two different conditions after it – 11 → 1 and 11 → 0. 11 → 1 (6 times), 11 → 0(1 time), 10 → 1(1 time) and 01 → 1 (1 time) => 6 + 1 + 1 – right predictions and 11 → one miss => 1/10 = 10%."

# --- Update sheet view: scroll so column B is left-most visible, select E9 ---
$ws.Range("E9").Select()
$excel.ActiveWindow.ScrollColumn = 2
